# The deck ships two themes:
#   theme1.xml -> the Slide Master's theme (was "Integral" / "Red Violet")
#   theme2.xml -> the Notes Master's theme (was "Office Theme" / "Office")
#
# The commit swaps the two themes' contents, so the Slide Master's theme
# becomes the stock "Office Theme" palette (and the Notes Master's theme
# becomes the "Integral" palette the Slide Master used to have).
#
# Recolor the Slide Master's ThemeColorScheme - the part of the Theme
# object model that actually drives the deck's on-screen appearance - to
# the swapped-in "Office Theme" palette.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme

# "Office Theme" color order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)

$tcs = $theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = ($hex -band 0xFF0000) -shr 16
    $g = ($hex -band 0x00FF00) -shr 8
    $b = ($hex -band 0x0000FF)
    $tcs.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
